$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 499.33334
$ws.Range("I41").Value = 499.33334
$ws.Range("K41").Value = 499.33334
$ws.Range("M41").Value = -59.33334000000002

$ws.Range("H62").Value = 4715.222
$ws.Range("I62").Value = 5996
$ws.Range("J62").Value = 4555.125
$ws.Range("K62").Value = 5996
$ws.Range("L62").Value = 4555.125
$ws.Range("M62").Value = -5372
$ws.Range("N62").Value = -5803.125

$ws.Range("H65").Value = 4715.222
$ws.Range("I65").Value = 5996
$ws.Range("J65").Value = 4555.125
$ws.Range("K65").Value = 29980
$ws.Range("L65").Value = 22775.625
$ws.Range("M65").Value = -26860
$ws.Range("N65").Value = -29015.625

$ws.Range("H111").Value = 5321.75
$ws.Range("I111").Value = 8469.75
$ws.Range("K111").Value = 25409.25
$ws.Range("M111").Value = -22342.25

$ws.Range("H116").Value = 11597.6
$ws.Range("I116").Value = 2996
$ws.Range("K116").Value = 2996
$ws.Range("M116").Value = 446

$ws.Range("H129").Value = 1715.5238
$ws.Range("J129").Value = 2616.4167
$ws.Range("L129").Value = 7849.250100000001
$ws.Range("N129").Value = -17849.2501

$ws.Range("H132").Value = 5101.8276
$ws.Range("I132").Value = 5509.731
$ws.Range("J132").Value = 1566.6666
$ws.Range("K132").Value = 16529.193
$ws.Range("L132").Value = 4699.9998
$ws.Range("M132").Value = -13999.193
$ws.Range("N132").Value = -9759.9998

$ws.Range("H137").Value = 6146.9395
$ws.Range("I137").Value = 2600.2307
$ws.Range("J137").Value = 19320.428
$ws.Range("K137").Value = 7800.6921
$ws.Range("L137").Value = 57961.284
$ws.Range("M137").Value = -5250.6921
$ws.Range("N137").Value = -63061.284

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2253.7
$ws.Range("I61").Value = 1497.1852
$ws.Range("J61").Value = 3824.923
$ws.Range("K61").Value = 1497.1852
$ws.Range("L61").Value = 3824.923
$ws.Range("M61").Value = -1285.1852
$ws.Range("N61").Value = -4248.923

$ws.Range("H74").Value = 160836.38
$ws.Range("I74").Value = 243237.61
$ws.Range("K74").Value = 243237.61
$ws.Range("M74").Value = -242363.61

$ws.Range("H77").Value = 160836.38
$ws.Range("I77").Value = 243237.61
$ws.Range("K77").Value = 1216188.05
$ws.Range("M77").Value = -1211820.05

$ws.Range("H122").Value = 3254
$ws.Range("I122").Value = 4017.4
$ws.Range("K122").Value = 12052.2
$ws.Range("M122").Value = -9602.200000000001

$ws.Range("H132").Value = 3390.3076
$ws.Range("I132").Value = 3341.5557
$ws.Range("K132").Value = 10024.6671
$ws.Range("M132").Value = -7494.667099999999

$ws.Range("H136").Value = 2253.7
$ws.Range("I136").Value = 1497.1852
$ws.Range("J136").Value = 3824.923
$ws.Range("K136").Value = 4491.5556
$ws.Range("L136").Value = 11474.769
$ws.Range("M136").Value = -1941.5556
$ws.Range("N136").Value = -16574.769

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H29").Value = 425
$ws.Range("I29").Value = 425
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 425
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = -136
$ws.Range("N29").ClearContents()

$ws.Range("H86").Value = 2621.7917
$ws.Range("I86").Value = 2010.1333
$ws.Range("J86").Value = 3641.2222
$ws.Range("K86").Value = 2010.1333
$ws.Range("L86").Value = 3641.2222
$ws.Range("M86").Value = -887.1333
$ws.Range("N86").Value = -5887.2222

$ws.Range("H89").Value = 2621.7917
$ws.Range("I89").Value = 2010.1333
$ws.Range("J89").Value = 3641.2222
$ws.Range("K89").Value = 10050.6665
$ws.Range("L89").Value = 18206.111
$ws.Range("M89").Value = -4434.666499999999
$ws.Range("N89").Value = -29438.111

$ws.Range("H99").Value = 93763.45
$ws.Range("I99").Value = 144985.72
$ws.Range("K99").Value = 144985.72
$ws.Range("M99").Value = -143487.72

$ws.Range("H107").Value = 2405230.2
$ws.Range("I107").Value = 2748580.2
$ws.Range("J107").Value = 1781.5
$ws.Range("K107").Value = 2748580.2
$ws.Range("L107").Value = 1781.5
$ws.Range("M107").Value = -2746660.2
$ws.Range("N107").Value = -5621.5

$ws.Range("H134").Value = 2036.7
$ws.Range("I134").Value = 1409.5714
$ws.Range("K134").Value = 4228.7142
$ws.Range("M134").Value = -1693.7142

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2781982.2
$ws.Range("I31").Value = 3293.7144
$ws.Range("J31").Value = 5213334.5
$ws.Range("K31").Value = 3293.7144
$ws.Range("L31").Value = 5213334.5
$ws.Range("M31").Value = -2998.7144
$ws.Range("N31").Value = -5213924.5

$ws.Range("H34").Value = 2781982.2
$ws.Range("I34").Value = 3293.7144
$ws.Range("J34").Value = 5213334.5
$ws.Range("K34").Value = 3293.7144
$ws.Range("L34").Value = 5213334.5
$ws.Range("M34").Value = -3091.7144
$ws.Range("N34").Value = -5213738.5

$ws.Range("H132").Value = 4379.515
$ws.Range("I132").Value = 4275
$ws.Range("J132").Value = 4540.3076
$ws.Range("K132").Value = 12825
$ws.Range("L132").Value = 13620.9228
$ws.Range("M132").Value = -10295
$ws.Range("N132").Value = -18680.9228

$ws.Range("H134").Value = 4476.852
$ws.Range("I134").Value = 5014.2
$ws.Range("K134").Value = 15042.6
$ws.Range("M134").Value = -12507.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1099
$ws.Range("J5").Value = 800
$ws.Range("L5").Value = 2400
$ws.Range("N5").Value = -2624

$ws.Range("H74").Value = 26872.75
$ws.Range("J74").Value = 27494.834
$ws.Range("L74").Value = 82484.50199999999
$ws.Range("N74").Value = -84606.50199999999

$ws.Range("H77").Value = 26872.75
$ws.Range("J77").Value = 27494.834
$ws.Range("L77").Value = 247453.506
$ws.Range("N77").Value = -258061.506

$ws.Range("H135").Value = 1099
$ws.Range("J135").Value = 800
$ws.Range("L135").Value = 7200
$ws.Range("N135").Value = -12270

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 5328.893
$ws.Range("I102").Value = 1456.1428
$ws.Range("J102").Value = 6619.8096
$ws.Range("K102").Value = 1456.1428
$ws.Range("L102").Value = 6619.8096
$ws.Range("M102").Value = 165.8571999999999
$ws.Range("N102").Value = -9863.809600000001

$ws.Range("H107").Value = 10382.8
$ws.Range("J107").Value = 14914.667
$ws.Range("L107").Value = 14914.667
$ws.Range("N107").Value = -18754.667

$ws.Range("H122").Value = 3209510.5
$ws.Range("I122").Value = 6996913
$ws.Range("K122").Value = 20990739
$ws.Range("M122").Value = -20988289

$ws.Range("H132").Value = 2164.2273
$ws.Range("I132").Value = 1658.1111
$ws.Range("K132").Value = 4974.3333
$ws.Range("M132").Value = -2444.3333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1497.8334
$ws.Range("J22").Value = 500
$ws.Range("L22").Value = 500
$ws.Range("N22").Value = -1090

$ws.Range("H27").Value = 1497.8334
$ws.Range("J27").Value = 500
$ws.Range("L27").Value = 500
$ws.Range("N27").Value = -714

$ws.Range("H61").Value = 1523.5
$ws.Range("I61").Value = 1523.5
$ws.Range("K61").Value = 1523.5
$ws.Range("M61").Value = -1321.5

$ws.Range("H93").Value = 2838.7646
$ws.Range("J93").Value = 2499.3333
$ws.Range("L93").Value = 2499.3333
$ws.Range("N93").Value = -4995.3333

$ws.Range("H113").Value = 1523.5
$ws.Range("I113").Value = 1523.5
$ws.Range("K113").Value = 1523.5
$ws.Range("M113").Value = 646.5

$ws.Range("H122").Value = 12997.6
$ws.Range("I122").Value = 8332.666999999999
$ws.Range("K122").Value = 24998.001
$ws.Range("M122").Value = -22548.001

$ws.Range("H132").Value = 6392.095
$ws.Range("I132").Value = 7078.5
$ws.Range("K132").Value = 21235.5
$ws.Range("M132").Value = -18705.5

$ws.Range("H137").Value = 43799.4
$ws.Range("I137").Value = 27599.8
$ws.Range("J137").Value = 59999
$ws.Range("K137").Value = 27599.8
$ws.Range("L137").Value = 59999
$ws.Range("M137").Value = -22499.8
$ws.Range("N137").Value = -70199

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()

$ws.Range("H132").Value = 1364.1945
$ws.Range("I132").Value = 1287.5358
$ws.Range("K132").Value = 3862.6074
$ws.Range("M132").Value = -1332.6074

$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()

$ws.Range("H136").Value = 4469.9287
$ws.Range("I136").Value = 2380.12
$ws.Range("J136").Value = 21885
$ws.Range("K136").Value = 7140.36
$ws.Range("L136").Value = 65655
$ws.Range("M136").Value = -4590.36
$ws.Range("N136").Value = -70755
